$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 382
$ws.Range("F7").Value = 1206
$ws.Range("F8").Value = 453
$ws.Range("F9").Value = 7397
$ws.Range("F11").Value = 95
$ws.Range("F13").Value = 8050
$ws.Range("F16").Value = 5537
$ws.Range("F18").Value = 2462
$ws.Range("F19").Value = 1047
$ws.Range("F25").Value = 408
$ws.Range("F26").Value = 405
$ws.Range("F28").Value = 2484
$ws.Range("F30").Value = 278
$ws.Range("F31").Value = 92
$ws.Range("F32").Value = 181
$ws.Range("F33").Value = 609
$ws.Range("F36").Value = 1550
$ws.Range("F38").Value = 10
$ws.Range("F39").Value = 2433
$ws.Range("F40").Value = 2231

$ws = $wb.Worksheets.Item("演出")
$ws.Range("G2").Value = 108
$ws.Range("F3").Value = 88
$ws.Range("F4").Value = 81
$ws.Range("F5").Value = 22

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("G3").Value = "不可售"

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("G4").Value = "不可售"
$ws.Range("G6").Value = 108
$ws.Range("F7").Value = 382
$ws.Range("F8").Value = 1206
$ws.Range("F9").Value = 453
$ws.Range("F10").Value = 7397
$ws.Range("F12").Value = 95
$ws.Range("F14").Value = 8050
$ws.Range("F17").Value = 5537
$ws.Range("F19").Value = 2462
$ws.Range("F20").Value = 1047
$ws.Range("F24").Value = 88
$ws.Range("F26").Value = 81
$ws.Range("F27").Value = 408
$ws.Range("F28").Value = 405
$ws.Range("F30").Value = 2484
$ws.Range("F32").Value = 278
$ws.Range("F33").Value = 92
$ws.Range("F34").Value = 181
$ws.Range("F35").Value = 22
$ws.Range("F36").Value = 609
$ws.Range("F40").Value = 1550
$ws.Range("F42").Value = 10
$ws.Range("F43").Value = 2433
$ws.Range("F45").Value = 2231
